$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- Drop the oldest listing (old row 2) so the remaining two shift up ---
$ws.Rows("2:2").Delete()

# --- Drop the two trailing listings that were scraped away (old rows 5:6,
#     now rows 4:5 after the shift above) ---
$ws.Rows("4:5").Delete()

# --- Refresh the "fetched at" timestamp on the two surviving rows ---
$ws.Range("A2").Value = "2025-12-27 06:28:03"
$ws.Range("A3").Value = "2025-12-27 06:28:03"

# --- Hyperlinks: the engine doesn't re-target hyperlinks on row delete,
#     so rebuild the collection to match the (now 2-row) data ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5462048")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5461891")
# Hyperlinks.Add() re-applies the "Hyperlink" cell style through a fresh
# style record; put both cells back on the original Hyperlink style so the
# workbook keeps using the same style index it already had.
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"

# --- Column width tweaks (B: 41 -> 34, H: 13 -> 12). ColumnWidth uses
#     Excel's character-width units which store with a ~0.8333 padding
#     offset, so back that out to land on the exact stored width. ---
$ws.Columns.Item(2).ColumnWidth = 34 - 5/6
$ws.Columns.Item(8).ColumnWidth = 12 - 5/6
